$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '39.954.62'
$ws.Cells.Item(2, 5).Value = '  -0.32%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.209.21'
$ws.Cells.Item(3, 5).Value = '  -1.28%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  -0.18%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''294.18'
$ws.Cells.Item(5, 5).Value = '  -0.19%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''87.06'
$ws.Cells.Item(6, 5).Value = '  +0.94%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.509'
$ws.Cells.Item(7, 5).Value = '  -1.17%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.11%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.473'
$ws.Cells.Item(9, 5).Value = '  +0.15%  '

# Row 10
$ws.Cells.Item(10, 2).Value = 'Avalanche'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(10, 4).Value = '''29.87'
$ws.Cells.Item(10, 5).Value = '  -4.13%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'Dogecoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(11, 4).Value = '''0.0776'
$ws.Cells.Item(11, 5).Value = '  -2.03%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''49.11'
$ws.Cells.Item(12, 5).Value = '  +4.45%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +2.72%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''6.48'
$ws.Cells.Item(14, 5).Value = '  +0.20%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.547.82'
$ws.Cells.Item(15, 5).Value = '  -1.28%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''13.69'
$ws.Cells.Item(16, 5).Value = '  -3.28%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.214.68'
$ws.Cells.Item(17, 5).Value = '  +1.31%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''0.727'
$ws.Cells.Item(18, 5).Value = '  -0.39%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '39.811.00'
$ws.Cells.Item(19, 5).Value = '  -0.52%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0884'
$ws.Cells.Item(20, 5).Value = '  -0.86%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''11.34'
$ws.Cells.Item(21, 5).Value = '  +4.82%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''5.77'
$ws.Cells.Item(22, 5).Value = '  -0.72%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''65.12'
$ws.Cells.Item(23, 5).Value = '  -0.44%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''236.08'
$ws.Cells.Item(24, 5).Value = '  +0.50%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.04%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.78%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''1.81'
$ws.Cells.Item(27, 5).Value = '  -2.19%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''22.51'
$ws.Cells.Item(28, 5).Value = '  -1.48%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -3.52%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''9.18'
$ws.Cells.Item(30, 5).Value = '  -0.58%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''155.36'
$ws.Cells.Item(31, 5).Value = '  +2.01%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''31.77'
$ws.Cells.Item(32, 5).Value = '  -4.81%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -0.11%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''4.90'
$ws.Cells.Item(34, 5).Value = '  +0.21%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.0713'
$ws.Cells.Item(35, 5).Value = '  -1.25%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''2.33'
$ws.Cells.Item(36, 5).Value = '  -2.27%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''2.83'
$ws.Cells.Item(37, 5).Value = '  +4.05%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.18%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''0.0977'
$ws.Cells.Item(39, 5).Value = '  -2.58%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''15.55'
$ws.Cells.Item(40, 5).Value = '  -4.93%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''1.67'
$ws.Cells.Item(41, 5).Value = '  -1.74%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '2.127.88'
$ws.Cells.Item(42, 5).Value = '  +3.66%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''3.73'
$ws.Cells.Item(43, 5).Value = '  -2.90%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''2.11'
$ws.Cells.Item(44, 5).Value = '  -5.91%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'VeChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(45, 4).Value = '''0.0266'
$ws.Cells.Item(45, 5).Value = '  -1.46%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).Value = '''17.75'
$ws.Cells.Item(46, 5).Value = '  +8.87%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''9.63'
$ws.Cells.Item(47, 5).Value = '  -3.86%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''2.66'
$ws.Cells.Item(48, 5).Value = '  +3.73%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '2.415.15'
$ws.Cells.Item(49, 5).Value = '  -1.47%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''1.47'
$ws.Cells.Item(50, 5).Value = '  +0.38%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.38%  '
